# Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.
#
# Regenerating the document with the newer OOXML writer re-touches the
# character run that holds the cached result of each "REF bookmark1"
# field (the text "a reference to bookmark1"), which keeps its bold
# formatting. Re-apply Bold to that cached field-result text for both
# occurrences of the field (the link before the bookmark and the link
# after the bookmark) so the formatting survives the regeneration.

$d = $word.ActiveDocument

$searchText = "a reference to bookmark1"

# First occurrence: "Test link before bookmark : <field>"
$range1 = $d.Content
$found1 = $range1.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if ($found1) {
    $range1.Font.Bold = 1
}

# Second occurrence: "Test link after bookmark : <field>"
$range2 = $d.Range($range1.End, $d.Content.End)
$found2 = $range2.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if ($found2) {
    $range2.Font.Bold = 1
}
